$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F4").Value = 4
$ws.Range("F8").Value = -2
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("F30").Value = 3
$ws.Range("F31").Value = 1
$ws.Range("F35").Value = 2
$ws.Range("F41").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = -5
$ws.Range("F55").Value = 4
$ws.Range("F58").Value = -1
$ws.Range("F60").Value = -2
$ws.Range("F62").Value = -2
$ws.Range("F63").Value = 0
$ws.Range("F67").Value = 2
$ws.Range("F69").Value = -3
$ws.Range("F71").Value = 2
$ws.Range("F74").Value = 0
$ws.Range("F78").Value = 0
$ws.Range("F84").Value = -2
